$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# " ... prints out the square of that number. For example, an input of 5
#   should result in the number 25 being printed. "
# becomes
# " ... says the square of that number. For example, an input of 5
#   should result in the cat saying 25. "
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    " as an input and then prints out the square of that number. For example, an input of 5 should result in the number 25 being printed. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " as an input and then says the square of that number. For example, an input of 5 should result in the cat saying 25. ",
    2)

# --- Edit 2 -----------------------------------------------------------
# "problem to print all of the squares of the numbers between one and ten."
# becomes
# "problem to say all of the squares of the numbers between one and ten."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "problem to print all of the squares of the numbers between one and ten.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "problem to say all of the squares of the numbers between one and ten.",
    2)
